$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 373, pushing existing rows 373:402 down to 374:403.
$ws.Rows.Item(373).Insert()

# Populate the newly inserted row 373 with the new weekly price observation.
$ws.Cells.Item(373, 1).Value = 3
$ws.Cells.Item(373, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(373, 3).Value = "Coquimbo"
$ws.Cells.Item(373, 4).Value = 44783
$ws.Cells.Item(373, 5).Value = 5
$ws.Cells.Item(373, 6).Value = 100112040
$ws.Cells.Item(373, 7).Value = "Cilantro"
$ws.Cells.Item(373, 8).Value = "Sin especificar"
$ws.Cells.Item(373, 9).Value = "Primera"
$ws.Cells.Item(373, 10).Value = 115
$ws.Cells.Item(373, 11).Value = 4500
$ws.Cells.Item(373, 12).Value = 5000
$ws.Cells.Item(373, 13).Value = 4761
$ws.Cells.Item(373, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(373, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(373, 16).Value = 1587
$ws.Cells.Item(373, 17).Value = 3
$ws.Cells.Item(373, 18).Value = "Hortaliza"
